# Adding more dates to special dates
# Appends new "special dates" rows (116-162) to Hoja1, matching the
# highlighted (fill-themed) formatting already used by rows 112-115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1: rows 116-138 -> Dia / Mes / Periodo(2019), fully populated ----
$g1A = @(1,7,25,14,18,19,21,1,12,3,16,24,1,20,7,19,14,31,4,11,7,8,25)
$g1B = @(1,1,3,4,4,4,4,5,5,6,6,6,7,7,8,8,10,10,11,11,12,12,12)
$g1StartRow = 116

# Copy the highlighted format (style used by rows 112:115) onto the new block
$ws.Range("A112:C112").Copy() | Out-Null
$ws.Range("A116:C138").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $g1A.Length; $i++) {
    $r = $g1StartRow + $i
    $ws.Cells.Item($r, 1).Value = $g1A[$i]
    $ws.Cells.Item($r, 2).Value = $g1B[$i]
    $ws.Cells.Item($r, 3).Value = 2019
}

# ---- Block 2: rows 139-144 -> Dia / Mes only (no Periodo), but row still
# keeps a memory of column C (spans 1:3) because C briefly held a value
# before being cleared, same as the source workbook. ----
$g2A = @(1,6,23,5,9,10)
$g2B = @(1,1,3,4,4,4)
$g2StartRow = 139

$ws.Range("A112:C112").Copy() | Out-Null
$ws.Range("A139:C144").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $g2A.Length; $i++) {
    $r = $g2StartRow + $i
    $ws.Cells.Item($r, 1).Value = $g2A[$i]
    $ws.Cells.Item($r, 2).Value = $g2B[$i]
    $ws.Cells.Item($r, 3).Value = 9999
}
$ws.Range("C139:C144").Clear() | Out-Null

# ---- Block 3: rows 145-162 -> Dia / Mes only, column C untouched (spans 1:2) ----
$g3A = @(12,1,10,25,15,21,22,29,20,7,17,12,31,2,16,7,8,25)
$g3B = @(4,5,5,5,6,6,6,6,7,8,8,9,9,11,11,12,12,12)
$g3StartRow = 145

$ws.Range("A112:B112").Copy() | Out-Null
$ws.Range("A145:B162").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $g3A.Length; $i++) {
    $r = $g3StartRow + $i
    $ws.Cells.Item($r, 1).Value = $g3A[$i]
    $ws.Cells.Item($r, 2).Value = $g3B[$i]
}

# Reflect the new selection like the saved workbook did (active cell just
# past the last populated row).
$ws.Range("B163").Select() | Out-Null

Write-Output "Added rows 116-162 (47 new special dates)"
